$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1, J1 - copy formatting (style) from H1, then set values
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I (I0) and J (IF), rows 2-24
$data = @(
    @(8, 8),
    @(8, 9),
    @(1, 5),
    @(1, 4),
    @(2, 9),
    @(1, 4),
    @(1, 4),
    @(8, 9),
    @(2, 5),
    @(8, 8),
    @(8, 9),
    @(7, 8),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(9, 9),
    @(6, 7),
    @(8, 8),
    @(11, 11),
    @(4, 5),
    @(6, 7),
    @(7, 7),
    @(9, 9)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
